# Sync attendance_reports: swap the order of the two names/emails listed
# in the "Recorded By" column (G) so that "dnasr281@gmail.com" appears
# first, e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
# and "admin@admin.com, dnasr281@gmail.com" -> "dnasr281@gmail.com, admin@admin.com".
# Cells that only contain a single value, or more than two comma-separated
# values, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $val = $cell.Value2

    if ($val -ne $null -and $val.GetType().Name -eq "String" -and $val.Contains(",")) {
        $parts = $val.Split(",")
        if ($parts.Length -eq 2) {
            $first = $parts[0].Trim()
            $second = $parts[1].Trim()
            if ($first -eq "dnasr281@gmail.com" -or $second -ne "dnasr281@gmail.com") {
                continue
            }
            $cell.Value = "$second, $first"
        }
    }
}
